$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 56.7
$ws.Range("N2").Value = 49.16024380385575

$ws.Range("D3").Value = 15.35
$ws.Range("E3").Value = 57.2
$ws.Range("F3").Value = 6.34
$ws.Range("K3").Value = 56.7
$ws.Range("N3").Value = 49.16024380385575
